# Update the "timestamp" column (H) values for data rows 2-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = 45767.431891083

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 8).Value = $newTimestamp
}
